$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43 (shifts existing rows 43..110 down to 44..111,
# preserving all their values/formatting - matches the diff where every row's
# data moved down by one and a brand-new record appears at row 43).
$ws.Rows(43).Insert()

# Populate the newly inserted row 43 with the new weekly record. Columns A-C,
# E-J, N, O, Q, R keep the same values as the rest of the (unchanged) dataset
# for this market/category, while D, K, L, M and P carry the new figures.
$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44797
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112001
$ws.Range("G43").Value = "Berenjena"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 100
$ws.Range("K43").Value = 13000
$ws.Range("L43").Value = 14000
$ws.Range("M43").Value = 13500
$ws.Range("N43").Value = "$/caja 60 unidades"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 225
$ws.Range("Q43").Value = 60
$ws.Range("R43").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the rest of
# column D.
$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat()
